# Add a new attendance-date column (U) for 2025-06-23, update the
# running "Total" counts in column S, and mark the new day absent (❌)
# for every student, mirroring the existing daily-attendance columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell U1: new date column header -------------------------------
# Leading apostrophe keeps "2025-06-23" as literal text instead of being
# auto-converted to a date serial number.
$ws.Range("U1").Value = "'2025-06-23"

# Copy the formatting of the neighboring header cell (T1: "Attendance %")
# onto U1 so it picks up the same bold/centered/bordered header style.
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)

# --- Row 2 (Abhishek Pathak): bump Total, mark new day absent -------------
$ws.Range("S2").Value = 16
$ws.Range("U2").Value = "❌"

# --- Row 3 (Shubham Pitekar): bump Total, mark new day absent -------------
$ws.Range("S3").Value = 16
$ws.Range("U3").Value = "❌"

"done"
